# Automatische test-sync: 2025-06-19 19:45:30
#
# Adds a new incoming mail-log entry ("Vragen over samenwerking" /
# "Samenwerking / Partnerverzoek") to the Logs sheet (row 42) and the
# corresponding dashboard tally row (row 12) on the Dashboard sheet,
# then extends the conditional formatting ranges and the bar chart's
# category/value series references to include the new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Logs sheet: append the new row (row 42)
# ---------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A42").Value = "Vragen over samenwerking"
$logs.Range("B42").Value = "mailmind.test@zohomail.eu"
$logs.Range("C42").Value = "Kunnen we samenwerken aan een nieuw project?"
$logs.Range("D42").Value = "Samenwerking / Partnerverzoek"
$logs.Range("F42").Value = "2025-06-19 19:45:25"
$logs.Range("G42").Value = "Nee"

# Extend the conditional formatting ranges to cover the new row
$logs.Range("D2:D41").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D42"))
$logs.Range("G2:G41").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G42"))

# ---------------------------------------------------------------
# 2. Dashboard sheet: append the new tally row (row 12)
# ---------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A12").Value = "Samenwerking / Partnerverzoek"
$dash.Range("B12").Value = 1

# ---------------------------------------------------------------
# 3. Update the chart's category/value series references so they
#    include the new Dashboard row (A2:A12 / B2:B12)
# ---------------------------------------------------------------
$chart = $dash.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!B1,'Dashboard'!`$A`$2:`$A`$12,'Dashboard'!`$B`$2:`$B`$12,1)"
